# Apply the updated cryptocurrency price/volume snapshot (GitHub Actions data refresh).
# Numeric-looking values in column D (e.g. "1.001", "235.15") are written with a
# leading apostrophe so Excel keeps them as literal text (matching the source
# workbook's inline-string cells) instead of silently re-parsing them as numbers
# and losing the original formatting/precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '25.283.98'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '1.657.39'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').Value = '''235.15'
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('D7').Value = '''0.4767'
$ws.Range('E7').Value = '  -2.27%  '
$ws.Range('D8').Value = '''0.2595'
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('D9').Value = '''0.06127'
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('D10').Value = '''0.07061'
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D11').Value = '1.658.31'
$ws.Range('E11').Value = '  -1.61%  '
$ws.Range('D12').Value = '''14.66'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '''0.5879'
$ws.Range('E13').Value = '  -6.76%  '
$ws.Range('D14').Value = '''4.369'
$ws.Range('D15').Value = '''74.18'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = '''1.001'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '''1.002'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '25.321.38'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').Value = '''0.000006726'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').Value = '''11.38'
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('D21').Value = '1.874.43'
$ws.Range('E21').Value = '  -2.25%  '
$ws.Range('D22').Value = '''4.418'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').Value = '''8.613'
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('D24').Value = '''5.313'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').Value = '''133.01'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('D26').Value = '''15.03'
$ws.Range('E26').Value = '  +0.85%  '
$ws.Range('D27').Value = '''1.398'
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('D28').Value = '''103.94'
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').Value = '''1.681'
$ws.Range('E29').Value = '  -2.79%  '
$ws.Range('D30').Value = '''3.970'
$ws.Range('E30').Value = '  +3.36%  '
$ws.Range('D31').Value = '''3.601'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('D32').Value = '''0.07619'
$ws.Range('E32').Value = '  -4.49%  '
$ws.Range('D33').Value = '''0.04350'
$ws.Range('E33').Value = '  -6.14%  '
$ws.Range('D34').Value = '''1.001'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('D36').Value = '''0.6085'
$ws.Range('E36').Value = '  +3.46%  '
$ws.Range('D37').Value = '''0.9400'
$ws.Range('E37').Value = '  -2.24%  '
$ws.Range('D38').Value = '''2.604'
$ws.Range('E38').Value = '  -2.74%  '
$ws.Range('D39').Value = '''0.8524'
$ws.Range('E39').Value = '  +0.89%  '
$ws.Range('D40').Value = '''1.000'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').Value = '''0.01498'
$ws.Range('E41').Value = '  -4.61%  '
$ws.Range('D42').Value = '''1.816'
$ws.Range('E42').Value = '  -4.48%  '
$ws.Range('D43').Value = '''97.79'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('D44').Value = '''0.3750'
$ws.Range('E44').Value = '  -0.79%  '
$ws.Range('D45').Value = '''4.641'
$ws.Range('E45').Value = '  -5.89%  '
$ws.Range('D46').Value = '''6.173'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '''0.1102'
$ws.Range('E47').Value = '  -4.37%  '
$ws.Range('D48').Value = '''0.05249'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('D49').Value = '''29.37'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').Value = '''1.002'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '''1.002'
$ws.Range('E51').Value = '  -0.41%  '
